# Calculate error of 2 BSSIDs with highest RSSI values
# (LocationAPI without RSSI in request) - fill in column C
# on the "Highest RSSI" sheet, and drop the now-unused column F
# placeholder cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Highest RSSI")

# --- Fill in the new "LocationAPI: error of combination with highest RSSI2"
#     column (C) with the computed error values / "No match" markers ---
$ws.Range("C2").Value = 0.0036777683199058619
$ws.Range("C3").Value = 0.056295732946324617
$ws.Range("C4").Value = 0.28018565116346822
$ws.Range("C5").Value = 0.09108448859608019
$ws.Range("C6").Value = 0.075915521684292817
$ws.Range("C7").Value = 0.0063388027275901159
$ws.Range("C8").Value = 0.021970540768343809
$ws.Range("C9").Value = 0.04882855601271921
$ws.Range("C11").Value = 0.010264361308076749
$ws.Range("C12").Value = 0.05850408241707817
$ws.Range("C13").Value = 0.02784041719545963
$ws.Range("C14").Value = 0.051323973692984859
$ws.Range("C15").Value = 0.0046433748571819301
$ws.Range("C16").Value = 0.044014910901996003
$ws.Range("C18").Value = 0.032644514686287388
$ws.Range("C19").Value = 0.026481671533126909
$ws.Range("C20").Value = 0.0034100794354838941
$ws.Range("C21").Value = 0.058719428050553413
$ws.Range("C23").Value = "No match"
$ws.Range("C24").Value = 0.02980599690868627
$ws.Range("C25").Value = "No match"
$ws.Range("C27").Value = 0.044046625330704393
$ws.Range("C28").Value = 0.036462036747490233
$ws.Range("C29").Value = 0.031613731843301768
$ws.Range("C31").Value = 0.01325813851648
$ws.Range("C32").Value = 0.04257717226305445
$ws.Range("C33").Value = "No match"
$ws.Range("C34").Value = "No match"
$ws.Range("C35").Value = "No match"
$ws.Range("C36").Value = "No match"
$ws.Range("C37").Value = "No match"

# --- Remove the stray, unused placeholder cells in column F (rows 2-13) ---
$ws.Range("F2:F13").Clear()

# --- Widen column C a touch so the new data fits nicely ---
$ws.Columns("C").ColumnWidth = 47.6

# --- Update the active selection to match the edited sheet state ---
$ws.Range("F6").Select()
